# Update countries & provincias Spain
# Applies the COVID data refresh captured in the target diff:
#  - Armenia overtakes Rumania in total cases -> rows 52/53 swap order/labels
#  - Groenlandia/Islas Malvinas swap display order (identical figures)
#  - Updated case counters for several countries (rows 4,6,7,38,41,52,53,97,114,119,132,192)
#  - "Datos actualizados" timestamp bumped from 08:31 to 09:48

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 7 de Julio de 2020 a las 09:48"

# --- Updated statistics for existing countries -------------------------
# Row 4: Estados Unidos
Set-Row 4 3041035 202 1325066 1582988 0 2 132981

# Row 6: India
Set-Row 6 721310 964 440229 260897 0 10 20184

# Row 7: Rusia
Set-Row 7 694230 6368 463880 219856 0 198 10494

# Row 38: Ucrania
Set-Row 38 49607 564 22193 26131 0 21 1283

# Row 41: Singapur
Set-Row 41 45140 157 40717 4397 0 0 26

# Rows 52/53: Armenia overtakes Rumania, so the two rows swap country and data
$ws.Range("A52").Value = "Armenia"
Set-Row 52 29285 349 16907 11875 0 12 503

$ws.Range("A53").Value = "Rumania"
Set-Row 53 29223 0 20213 7242 0 0 1768

# Row 97: Hungria
Set-Row 97 4205 16 2874 742 0 0 589

# Row 114: Estonia
Set-Row 114 1995 1 1880 46 0 0 69

# Row 119: Eslovaquia
Set-Row 119 1767 2 1473 266 0 0 28

# Row 132: Letonia
Set-Row 132 1134 7 1008 96 0 0 30

# Row 192: Islas Turcas y Caicos
Set-Row 192 49 1 11 36 0 0 2

# Rows 209/210: Groenlandia and Islas Malvinas swap display order (figures unchanged)
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"
